$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item(1)

# --- Create the new "TrackingOrder" sheet right after "Login" ---
$tracking = $wb.Worksheets.Add([System.Type]::Missing, $login)
$tracking.Name = "TrackingOrder"

# Header typed first (this is why "OrderId" is the first *new* shared string)
$tracking.Range("A1").Value = "OrderId"

# --- Now fill in the Login sheet ---
# NOTE: NumberFormat is applied *before* Value everywhere a Text format is
# used, so that numeric-looking strings (leading zeros, etc.) are not
# silently coerced into numbers.

# Row 2
$login.Range("B2").NumberFormat = "@"
$login.Range("B2").Value = "Hatemyself@1001@@"
$login.Hyperlinks.Add($login.Range("B2"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 3
$login.Range("A3").NumberFormat = "@"
$login.Range("A3").Value = "thuctanphu12@gmail.com"
$login.Hyperlinks.Add($login.Range("A3"), "mailto:thuctanphu12@gmail.com") | Out-Null

$login.Range("B3").NumberFormat = "@"
$login.Range("B3").Value = "123456"

# Row 4
$login.Range("A4").NumberFormat = "@"
$login.Range("A4").Value = "thuctanphu12@gmail.com"
$login.Hyperlinks.Add($login.Range("A4"), "mailto:thuctanphu12@gmail.com") | Out-Null

$login.Range("B4").NumberFormat = "@"
$login.Range("B4").Value = "Hatemyself@1001@@"
$login.Hyperlinks.Add($login.Range("B4"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 5 (hyperlink only, left in General format)
$login.Range("A5").Value = "aaa@bbb"
$login.Hyperlinks.Add($login.Range("A5"), "mailto:aaa@bbb") | Out-Null

$login.Range("B5").NumberFormat = "@"
$login.Range("B5").Value = "Hatemyself@1001@@"
$login.Hyperlinks.Add($login.Range("B5"), "mailto:Hatemyself@1001@@") | Out-Null

# A2 is (re)entered at this point in the original authoring sequence
$login.Range("A2").NumberFormat = "@"
$login.Range("A2").Value = "03547305"

# Row 6
$login.Range("A6").NumberFormat = "@"
$login.Range("A6").Value = "thuctanphu12@gmail.com"
$login.Hyperlinks.Add($login.Range("A6"), "mailto:thuctanphu12@gmail.com") | Out-Null

# Row 7 (only column B)
$login.Range("B7").NumberFormat = "@"
$login.Range("B7").Value = "Hatemyself@1001@@"
$login.Hyperlinks.Add($login.Range("B7"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 8
$login.Range("A8").Value = 354730579

$login.Range("B8").NumberFormat = "@"
$login.Range("B8").Value = "Hatemyself@1001@@"
$login.Hyperlinks.Add($login.Range("B8"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 9
$login.Range("A9").Value = "thuctanphuaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaa@gmail.com"

$login.Range("B9").NumberFormat = "@"
$login.Range("B9").Value = "Hatemyself@1001@@"
$login.Hyperlinks.Add($login.Range("B9"), "mailto:Hatemyself@1001@@") | Out-Null

# Column A width on Login
$login.Columns.Item(1).ColumnWidth = 8.93

# --- Finish the TrackingOrder sheet ---
$tracking.Range("A2").Value = 392921444522425
$tracking.Range("A3").Value = 440277016458018
$tracking.Range("A4").Value = 1
$tracking.Range("A5").Value = "123456789a"

# --- Selections: Login ends up scrolled to B9, TrackingOrder is the active sheet at G14 ---
$login.Application.Goto($login.Range("B9"))
$tracking.Application.Goto($tracking.Range("G14"))
